$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4672717750072479
$ws.Range("B1").Value = 0.79325270652771
$ws.Range("C1").Value = 5.575876712799072
$ws.Range("D1").Value = 1.567294716835022
$ws.Range("E1").Value = 0.8993276953697205
